$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- Step 1: insert two new paragraphs right after the first paragraph ---
# (the first paragraph currently reads "Fffff...ttttt" and stays untouched)
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()

# New paragraph 2: spellStart + two runs ("F" / "ffff...") + spellEnd
$p2 = $d.Paragraphs(2)
$xmlP2 = '<w:p xmlns:w="' + $wNs + '">' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:t>F</w:t></w:r>' +
           '<w:r><w:t>ffffffffffffffffffffffffffffffffffffffffffffffffffffffff</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
         '</w:p>'
$p2.Range.InsertXML($xmlP2)

# New paragraph 3: spellStart + single run (duplicate of paragraph 1's text) + spellEnd
$p3 = $d.Paragraphs(3)
$xmlP3 = '<w:p xmlns:w="' + $wNs + '">' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:t>Fffffffffffffffffffffffffffftttttttttttttttttttttt</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
         '</w:p>'
$p3.Range.InsertXML($xmlP3)

# --- Step 2: split the bookmark out of the last "ffff" paragraph into its own paragraph ---
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()

# Rewrite paragraph 4: wrap its run with spellStart/gramStart .. spellEnd/gramEnd, drop the bookmark
$p4 = $d.Paragraphs(4)
$xmlP4 = '<w:p xmlns:w="' + $wNs + '">' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:proofErr w:type="gramStart"/>' +
           '<w:r><w:t>fffffffffffffffffffffffffffffffffffffffffffffffffffffffff</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:proofErr w:type="gramEnd"/>' +
         '</w:p>'
$p4.Range.InsertXML($xmlP4)

# New paragraph 5: just the bookmark
$p5 = $d.Paragraphs(5)
$xmlP5 = '<w:p xmlns:w="' + $wNs + '">' +
           '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
           '<w:bookmarkEnd w:id="0"/>' +
         '</w:p>'
$p5.Range.InsertXML($xmlP5)
